# Auto-generated cell updates for Hades_Profits leve-profit sheets
# Applies per-cell numeric value corrections across ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR sheets
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H19").Value = 470.77274
$ws.Range("I19").Value = 315
$ws.Range("J19").Value = 559.7857
$ws.Range("K19").Value = 315
$ws.Range("L19").Value = 559.7857
$ws.Range("M19").Value = -140
$ws.Range("N19").Value = -909.7857
$ws.Range("H64").Value = 4775.533
$ws.Range("J64").Value = 4775.533
$ws.Range("L64").Value = 4775.533
$ws.Range("N64").Value = -5271.533
$ws.Range("H67").Value = 4775.533
$ws.Range("J67").Value = 4775.533
$ws.Range("L67").Value = 4775.533
$ws.Range("N67").Value = -6491.533
$ws.Range("H98").Value = 606.7406999999999
$ws.Range("I98").Value = 627.84
$ws.Range("K98").Value = 627.84
$ws.Range("M98").Value = 870.16
$ws.Range("H112").Value = 41668572
$ws.Range("J112").Value = 2227.7778
$ws.Range("L112").Value = 6683.3334
$ws.Range("N112").Value = -8899.3334
$ws.Range("H122").Value = 606.7406999999999
$ws.Range("I122").Value = 627.84
$ws.Range("K122").Value = 1883.52
$ws.Range("M122").Value = 566.48
$ws.Range("H138").Value = 4880466
$ws.Range("J138").Value = 6062879.5
$ws.Range("L138").Value = 18188638.5
$ws.Range("N138").Value = -18198918.5
$ws.Range("H140").Value = 54740.332
$ws.Range("J140").Value = 54740.332
$ws.Range("L140").Value = 54740.332
$ws.Range("N140").Value = -65100.332

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 9532.66
$ws.Range("I32").Value = 6101.317
$ws.Range("J32").Value = 25164.334
$ws.Range("K32").Value = 6101.317
$ws.Range("L32").Value = 25164.334
$ws.Range("M32").Value = -5814.317
$ws.Range("N32").Value = -25738.334
$ws.Range("H122").Value = 2116.6155
$ws.Range("I122").Value = 1650.25
$ws.Range("J122").Value = 2862.8
$ws.Range("K122").Value = 4950.75
$ws.Range("L122").Value = 8588.400000000001
$ws.Range("M122").Value = -2500.75
$ws.Range("N122").Value = -13488.4
$ws.Range("H132").Value = 45112.938
$ws.Range("I132").Value = 30840.656
$ws.Range("J132").Value = 80793.64
$ws.Range("K132").Value = 92521.96799999999
$ws.Range("L132").Value = 242380.92
$ws.Range("M132").Value = -89991.96799999999
$ws.Range("N132").Value = -247440.92

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H22").Value = 518
$ws.Range("J22").Value = 0
$ws.Range("L22").Value = 0
$ws.Range("N22").ClearContents()
$ws.Range("H94").Value = 1032.8
$ws.Range("J94").Value = 1010
$ws.Range("L94").Value = 1010
$ws.Range("N94").Value = -1912
$ws.Range("H105").Value = 2545.5557
$ws.Range("I105").Value = 2488.75
$ws.Range("K105").Value = 2488.75
$ws.Range("M105").Value = -741.75
$ws.Range("H107").Value = 3363.9
$ws.Range("I107").Value = 2543.2632
$ws.Range("J107").Value = 4781.364
$ws.Range("K107").Value = 2543.2632
$ws.Range("L107").Value = 4781.364
$ws.Range("M107").Value = -623.2631999999999
$ws.Range("N107").Value = -8621.364
$ws.Range("H141").Value = 58895
$ws.Range("J141").Value = 58895
$ws.Range("L141").Value = 58895
$ws.Range("N141").Value = -69255

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 271225.8
$ws.Range("I31").Value = 65132.5
$ws.Range("J31").Value = 384932.47
$ws.Range("K31").Value = 65132.5
$ws.Range("L31").Value = 384932.47
$ws.Range("M31").Value = -64837.5
$ws.Range("N31").Value = -385522.47
$ws.Range("H34").Value = 271225.8
$ws.Range("I34").Value = 65132.5
$ws.Range("J34").Value = 384932.47
$ws.Range("K34").Value = 65132.5
$ws.Range("L34").Value = 384932.47
$ws.Range("M34").Value = -64930.5
$ws.Range("N34").Value = -385336.47
$ws.Range("H52").Value = 38754.5
$ws.Range("J52").Value = 46800
$ws.Range("L52").Value = 46800
$ws.Range("N52").Value = -47388
$ws.Range("H99").Value = 1800
$ws.Range("I99").Value = 1800
$ws.Range("K99").Value = 1800
$ws.Range("M99").Value = -302
$ws.Range("H126").Value = 1800
$ws.Range("I126").Value = 1800
$ws.Range("K126").Value = 5400
$ws.Range("M126").Value = -2930
$ws.Range("H132").Value = 31568.883
$ws.Range("I132").Value = 1886.3077
$ws.Range("J132").Value = 128037.25
$ws.Range("K132").Value = 5658.9231
$ws.Range("L132").Value = 384111.75
$ws.Range("M132").Value = -3128.9231
$ws.Range("N132").Value = -389171.75
$ws.Range("H137").Value = 40316.5
$ws.Range("J137").Value = 40316.5
$ws.Range("L137").Value = 40316.5
$ws.Range("N137").Value = -50516.5

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H86").Value = 696.8182
$ws.Range("I86").Value = 300
$ws.Range("J86").Value = 736.5
$ws.Range("K86").Value = 900
$ws.Range("L86").Value = 2209.5
$ws.Range("M86").Value = 286
$ws.Range("N86").Value = -4581.5
$ws.Range("H89").Value = 696.8182
$ws.Range("I89").Value = 300
$ws.Range("J89").Value = 736.5
$ws.Range("K89").Value = 2700
$ws.Range("L89").Value = 6628.5
$ws.Range("M89").Value = 3228
$ws.Range("N89").Value = -18484.5

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H102").Value = 1811.5714
$ws.Range("I102").Value = 1866.8572
$ws.Range("J102").Value = 1756.2858
$ws.Range("K102").Value = 1866.8572
$ws.Range("L102").Value = 1756.2858
$ws.Range("M102").Value = -244.8571999999999
$ws.Range("N102").Value = -5000.2858
$ws.Range("H122").Value = 3222.2222
$ws.Range("I122").Value = 2620
$ws.Range("J122").Value = 3975
$ws.Range("K122").Value = 7860
$ws.Range("L122").Value = 11925
$ws.Range("M122").Value = -5410
$ws.Range("N122").Value = -16825
$ws.Range("H126").Value = 1646.2858
$ws.Range("J126").Value = 2000
$ws.Range("L126").Value = 6000
$ws.Range("N126").Value = -10940

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H61").Value = 2265.6316
$ws.Range("I61").Value = 2196.375
$ws.Range("J61").Value = 2635
$ws.Range("K61").Value = 2196.375
$ws.Range("L61").Value = 2635
$ws.Range("M61").Value = -1994.375
$ws.Range("N61").Value = -3039
$ws.Range("H68").Value = 1540.1333
$ws.Range("I68").Value = 1370.2
$ws.Range("J68").Value = 1880
$ws.Range("K68").Value = 1370.2
$ws.Range("L68").Value = 1880
$ws.Range("M68").Value = -621.2
$ws.Range("N68").Value = -3378
$ws.Range("H71").Value = 1540.1333
$ws.Range("I71").Value = 1370.2
$ws.Range("J71").Value = 1880
$ws.Range("K71").Value = 6851
$ws.Range("L71").Value = 9400
$ws.Range("M71").Value = -3107
$ws.Range("N71").Value = -16888
$ws.Range("H82").Value = 1527.5
$ws.Range("I82").Value = 1285
$ws.Range("K82").Value = 1285
$ws.Range("M82").Value = -924
$ws.Range("H85").Value = 1527.5
$ws.Range("I85").Value = 1285
$ws.Range("K85").Value = 1285
$ws.Range("M85").Value = -37
$ws.Range("H100").Value = 1988.3
$ws.Range("I100").Value = 1975.75
$ws.Range("J100").Value = 1996.6666
$ws.Range("K100").Value = 1975.75
$ws.Range("L100").Value = 1996.6666
$ws.Range("M100").Value = -1434.75
$ws.Range("N100").Value = -3078.6666
$ws.Range("H113").Value = 2265.6316
$ws.Range("I113").Value = 2196.375
$ws.Range("J113").Value = 2635
$ws.Range("K113").Value = 2196.375
$ws.Range("L113").Value = 2635
$ws.Range("M113").Value = -26.375
$ws.Range("N113").Value = -6975

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H107").Value = 437.19446
$ws.Range("I107").Value = 376.69232
$ws.Range("J107").Value = 594.5
$ws.Range("K107").Value = 1130.07696
$ws.Range("L107").Value = 1783.5
$ws.Range("M107").Value = 789.9230400000001
$ws.Range("N107").Value = -5623.5
$ws.Range("H113").Value = 674.9394
$ws.Range("I113").Value = 878.7222
$ws.Range("J113").Value = 430.4
$ws.Range("K113").Value = 2636.1666
$ws.Range("L113").Value = 1291.2
$ws.Range("M113").Value = -466.1666
$ws.Range("N113").Value = -5631.2
